$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "efd1bbba1f2029f684155f82f3f1d86f"  # 05-050207TP: 855d5b9e061f03507213163ccf594c50 -> efd1bbba1f2029f684155f82f3f1d86f
$ws.Range("B24").Value = "b53816fd5416d69f8d08e64f33be22ea"  # 05-050316TC: 8878748c700d7fdea464e30c7e4067fb -> b53816fd5416d69f8d08e64f33be22ea
$ws.Range("B34").Value = "b996e797ba6212dd6aef9edede692abe"  # 05-050316TP: f26b8661c6953e87a1e135d1ce4ba4f9 -> b996e797ba6212dd6aef9edede692abe
$ws.Range("B123").Value = "b0bda9e33e30da7e47e182fd0cfd97b0"  # 05-050301TC: cf0d998a1b7e6de4a284e3d22e487f5e -> b0bda9e33e30da7e47e182fd0cfd97b0
$ws.Range("B133").Value = "1a84eba233200c1095bdf0ce05c69593"  # 05-050312TP: 4409143d57b5150097d5d36c17aa15f5 -> 1a84eba233200c1095bdf0ce05c69593
$ws.Range("B163").Value = "99b699eaa1a805bde9ca6d1d096f8362"  # 05-050308A: d02109d78d059c69e67930e83c9ddf01 -> 99b699eaa1a805bde9ca6d1d096f8362
$ws.Range("B181").Value = "ea99afab03ae848972b286d07f656d8e"  # 05-050303TC: 803a55a9f4255f6dc2a4d211ac6630fd -> ea99afab03ae848972b286d07f656d8e
$ws.Range("B192").Value = "07495a158447b32746484e329b4d6f69"  # 05-050314TP: 3969bb9ea333d1f3d19157823fe04a57 -> 07495a158447b32746484e329b4d6f69
$ws.Range("B199").Value = "663afb40a703acec0708b0eb2ace8f2b"  # 05-050314TC: e2cd9281650b2561cce6e981c5071842 -> 663afb40a703acec0708b0eb2ace8f2b
$ws.Range("B214").Value = "218218e3cbd37f47844084d863cf6539"  # 05-050303A: d0871aa81a228cdf44e3daa7125e033b -> 218218e3cbd37f47844084d863cf6539
$ws.Range("B426").Value = "930e9bd628ccd09c643cd2b4a4b8cfad"  # 05-0709-070905BTC: 0841f66eec1f7caf51680bed6f5054c6 -> 930e9bd628ccd09c643cd2b4a4b8cfad
$ws.Range("B491").Value = "dfaa531734479ff24c0cc86be34eeb26"  # 05-050314A: e14fe01c910387baaad5cba5ac23c98e -> dfaa531734479ff24c0cc86be34eeb26
$ws.Range("B515").Value = "ac0d09498744214d108e07d6bfb29fcf"  # 05-050208TP: d610ef912cbfe99c4c2415100db28a0d -> ac0d09498744214d108e07d6bfb29fcf
$ws.Range("B520").Value = "683ec326156b4727b51a8147b53d0579"  # 05-050306TP: 4675c67bf2a8dc16775ec05abb7d5af3 -> 683ec326156b4727b51a8147b53d0579
$ws.Range("B528").Value = "a1d1e120bc126e136cd517ceabe6c06c"  # 05-050317TC: 0c6ae3d60f548d9aaf15ba28c7ac77f9 -> a1d1e120bc126e136cd517ceabe6c06c
$ws.Range("B529").Value = "45d2bc6ca943aa88d293e4dd42e0cc8b"  # 05-050312A: 7febf5349f4310f03801db71426221db -> 45d2bc6ca943aa88d293e4dd42e0cc8b
$ws.Range("B539").Value = "99b92e08a9d9b95cd2d17d26064138e1"  # 05-050317TP: 927fdd666ff5c2131184c7611ca11117 -> 99b92e08a9d9b95cd2d17d26064138e1
$ws.Range("B651").Value = "07e836543dda5ddf8c90ab77ee7875a4"  # 05-050302TC: cff54a9e4d0702d29363119765df9c28 -> 07e836543dda5ddf8c90ab77ee7875a4
$ws.Range("B682").Value = "093d715d6ac8a954951ae9e1a15cb4c0"  # 05-050317A: 50d6b3928ae51952c9f11b33b97961e0 -> 093d715d6ac8a954951ae9e1a15cb4c0
$ws.Range("B720").Value = "80bb07e1be7c565344f735b929b5c8ed"  # 05-050315A: 194d96116d0b83bc7346b5f030d7ef73 -> 80bb07e1be7c565344f735b929b5c8ed
$ws.Range("B742").Value = "fe1ab0e4cc668b481b2c83d103d09350"  # 05-050315TC: c3d5a10641f32913d3775147256cad50 -> fe1ab0e4cc668b481b2c83d103d09350
$ws.Range("B745").Value = "c28610927522e7d86b7c39af57e13235"  # 05-050316A: 0867eed9183bdebf6cc2ae2672c200c2 -> c28610927522e7d86b7c39af57e13235
$ws.Range("B749").Value = "869381debc71a2aaec56a9f2364eab30"  # 05-050207A: b75d3247d26a2130f844dc55796296cb -> 869381debc71a2aaec56a9f2364eab30
$ws.Range("B758").Value = "442f1d0c4688f48a08a7d611f758499f"  # 05-050315TP: a5326aa5e29f014ac571870f665bb88d -> 442f1d0c4688f48a08a7d611f758499f
$ws.Range("B872").Value = "e93fbf982867c3b48547b1b8085b9879"  # 05-050309TC: d8debcc55b9615a4dd4da9181ecdba57 -> e93fbf982867c3b48547b1b8085b9879
$ws.Range("B892").Value = "c9c849f03081bb7a17b5eba5feebb7ea"  # 03-030032A: d878f735a89572d2273c1e98708e28dd -> c9c849f03081bb7a17b5eba5feebb7ea
